$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Copy header + row formatting from existing table (rows 2 and 3-9) ----
$ws.Range("B2:F2").Copy() | Out-Null
$ws.Range("B11:F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$ws.Range("B3:F9").Copy() | Out-Null
$ws.Range("B12:F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F3:F9").Copy() | Out-Null
$ws.Range("G12:G23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Plain-value cells ----
$ws.Range("B11").Value = "Phase"
$ws.Range("C11").Value = "Current Status"
$ws.Range("D11").Value = "Trigger Event"
$ws.Range("E11").Value = "Next Action / Table"
$ws.Range("F11").Value = "Snowflake SQL Actions"
$ws.Range("G11").Value = "SP Call Example"
$ws.Range("B12").Value = "Initial Intake"
$ws.Range("C12").Value = "— (new claim)"
$ws.Range("D12").Value = "Claim ingestion"
$ws.Range("E12").Value = "➔ Insert into CLAIMS ➔ Insert into CLAIMS_QUEUE as OUTSTANDING ➔ Log in CLAIMS_QUEUE_AUDIT"
$ws.Range("G12").Value = "CALL proc_ingest_claim('CLAIM001', 'api_user', 'initial load')"
$ws.Range("B13").Value = "CLAIMS_QUEUE"
$ws.Range("C13").Value = "OUTSTANDING"
$ws.Range("D13").Value = "Accept"
$ws.Range("E13").Value = "➔ Remove from CLAIMS_QUEUE ➔ Insert into CLAIMS_MANAGER (NEW) ➔ Update CLAIMS to ACCEPTED ➔ Log in CLAIMS_QUEUE_AUDIT"
$ws.Range("G13").Value = "CALL proc_accept_claims('''CLAIM001'', ''CLAIM002''', 'api_user', 'bulk accept')"
$ws.Range("B14").Value = "CLAIMS_QUEUE"
$ws.Range("C14").Value = "OUTSTANDING"
$ws.Range("D14").Value = "Reject"
$ws.Range("E14").Value = "➔ Remove from CLAIMS_QUEUE ➔ Update CLAIMS to REJECTED ➔ Log in CLAIMS_QUEUE_AUDIT"
$ws.Range("G14").Value = "CALL proc_reject_claims('''CLAIM001'', ''CLAIM002''', 'api_user', 'bulk reject')"
$ws.Range("B15").Value = "CLAIMS_QUEUE"
$ws.Range("C15").Value = "REJECTED (by mistake)"
$ws.Range("D15").Value = "Reinstate"
$ws.Range("E15").Value = "➔ Insert back into CLAIMS_QUEUE as OUTSTANDING ➔ Update CLAIMS to OUTSTANDING ➔ Log in CLAIMS_QUEUE_AUDIT"
$ws.Range("G15").Value = "CALL proc_reinstate_claims('''CLAIM003''', 'api_user', 'reinstate rejected')"
$ws.Range("B16").Value = "CLAIMS_MANAGER"
$ws.Range("C16").Value = "NEW"
$ws.Range("D16").Value = "Investigator starts investigation"
$ws.Range("E16").Value = "➔ Update CLAIMS_MANAGER to UNDER_INVESTIGATION ➔ Log in CLAIMS_MANAGER_AUDIT"
$ws.Range("F16").Value = "UPDATE CLAIMS_MANAGER SET STATUS = 'UNDER_INVESTIGATION' WHERE CLAIM_ID = :claim_id"
$ws.Range("G16").Value = "CALL proc_start_investigation('''CLAIM004''', 'investigator_user', 'investigation opened')"
$ws.Range("B17").Value = "CLAIMS_MANAGER"
$ws.Range("C17").Value = "UNDER_INVESTIGATION"
$ws.Range("D17").Value = "Investigation resolved (approved)"
$ws.Range("E17").Value = "➔ Update CLAIMS_MANAGER to CLOSED ➔ Update CLAIMS to CLOSED ➔ Log in CLAIMS_MANAGER_AUDIT"
$ws.Range("G17").Value = "CALL proc_close_claims('''CLAIM004''', 'investigator_user', 'investigation approved')"
$ws.Range("B18").Value = "CLAIMS_MANAGER"
$ws.Range("C18").Value = "UNDER_INVESTIGATION"
$ws.Range("D18").Value = "Investigation resolved (rejected)"
$ws.Range("E18").Value = "➔ Update CLAIMS_MANAGER to REJECTED ➔ Update CLAIMS to REJECTED ➔ Log in CLAIMS_MANAGER_AUDIT"
$ws.Range("G18").Value = "CALL proc_manager_reject_claims('''CLAIM005''', 'investigator_user', 'investigation rejected')"
$ws.Range("B19").Value = "Any Phase"
$ws.Range("C19").Value = "Any"
$ws.Range("D19").Value = "Cancelation (by user request)"
$ws.Range("E19").Value = "➔ Remove from CLAIMS_QUEUE or CLAIMS_MANAGER ➔ Update CLAIMS to CANCELED ➔ Log in AUDIT"
$ws.Range("G19").Value = "CALL proc_cancel_claims('''CLAIM006''', 'user_requestor', 'cancel requested')"
$ws.Range("B20").Value = "CLAIMS_MANAGER"
$ws.Range("C20").Value = "CLOSED"
$ws.Range("D20").Value = "Reopen closed claim (appeal)"
$ws.Range("E20").Value = "➔ Update CLAIMS_MANAGER to UNDER_INVESTIGATION ➔ Update CLAIMS to UNDER_INVESTIGATION ➔ Log in CLAIMS_MANAGER_AUDIT"
$ws.Range("G20").Value = "CALL proc_reopen_claims('''CLAIM007''', 'appeal_user', 'reopened on appeal')"
$ws.Range("B21").Value = "Any Phase"
$ws.Range("C21").Value = "Any"
$ws.Range("D21").Value = "Bulk admin action (admin batch)"
$ws.Range("E21").Value = "➔ Apply bulk update or delete ➔ Update CLAIMS or QUEUE or MANAGER ➔ Log bulk change in AUDIT table"
$ws.Range("G21").Value = "CALL proc_admin_bulk_action('BULK_ACTION_TYPE', '''CLAIM008'', ''CLAIM009''', 'admin_user', 'bulk admin op')"
$ws.Range("B22").Value = "Any Phase"
$ws.Range("C22").Value = "Any"
$ws.Range("D22").Value = "Audit-only correction"
$ws.Range("E22").Value = "➔ Insert correction record into AUDIT tables only (no status change)"
$ws.Range("F22").Value = "INSERT INTO CLAIMS_QUEUE_AUDIT OR CLAIMS_MANAGER_AUDIT (...)"
$ws.Range("G22").Value = "CALL proc_audit_only_entry('''CLAIM010''', 'audit_user', 'manual correction entry')"
$ws.Range("B23").Value = "Any Phase"
$ws.Range("C23").Value = "OUTSTANDING or UNDER_INVESTIGATION"
$ws.Range("D23").Value = "Timed auto-close or auto-reject (timeout)"
$ws.Range("E23").Value = "➔ Check claim age ➔ Update CLAIMS or CLAIMS_MANAGER ➔ Log in appropriate AUDIT table"
$ws.Range("G23").Value = "CALL proc_auto_close_claims(:threshold_days, 'system_scheduler', 'timed auto-close/reject')"

# ---- Rich-text cells (column F, multiple runs; Calibri 11 default + Arial Unicode MS 10 segments) ----
$ws.Range("F12").Value = "INSERT INTO CLAIMS (CLAIM_ID, STATUS, CREATED_AT) VALUES (:claim_id, 'OUTSTANDING') INSERT INTO CLAIMS_QUEUE (...)"
$ws.Range("F12").Characters(85, 30).Font.Size = 10
$ws.Range("F12").Characters(85, 30).Font.Name = "Arial Unicode MS"

$ws.Range("F13").Value = "DELETE FROM CLAIMS_QUEUE WHERE CLAIM_ID = :claim_id INSERT INTO CLAIMS_MANAGER (...) UPDATE CLAIMS SET STATUS = 'ACCEPTED'"
$ws.Range("F13").Characters(53, 32).Font.Size = 10
$ws.Range("F13").Characters(53, 32).Font.Name = "Arial Unicode MS"
$ws.Range("F13").Characters(86, 37).Font.Size = 10
$ws.Range("F13").Characters(86, 37).Font.Name = "Arial Unicode MS"

$ws.Range("F14").Value = "DELETE FROM CLAIMS_QUEUE WHERE CLAIM_ID = :claim_id UPDATE CLAIMS SET STATUS = 'REJECTED'"
$ws.Range("F14").Characters(53, 37).Font.Size = 10
$ws.Range("F14").Characters(53, 37).Font.Name = "Arial Unicode MS"

$ws.Range("F15").Value = "INSERT INTO CLAIMS_QUEUE (CLAIM_ID, STATUS, REINSTATED_AT) VALUES (:claim_id, 'OUTSTANDING') UPDATE CLAIMS SET STATUS = 'OUTSTANDING'"
$ws.Range("F15").Characters(94, 40).Font.Size = 10
$ws.Range("F15").Characters(94, 40).Font.Name = "Arial Unicode MS"

$ws.Range("F17").Value = "UPDATE CLAIMS_MANAGER SET STATUS = 'CLOSED' WHERE CLAIM_ID = :claim_id UPDATE CLAIMS SET STATUS = 'CLOSED'"
$ws.Range("F17").Characters(72, 35).Font.Size = 10
$ws.Range("F17").Characters(72, 35).Font.Name = "Arial Unicode MS"

$ws.Range("F18").Value = "UPDATE CLAIMS_MANAGER SET STATUS = 'REJECTED' WHERE CLAIM_ID = :claim_id UPDATE CLAIMS SET STATUS = 'REJECTED'"
$ws.Range("F18").Characters(74, 37).Font.Size = 10
$ws.Range("F18").Characters(74, 37).Font.Name = "Arial Unicode MS"

$ws.Range("F19").Value = "DELETE FROM CLAIMS_QUEUE OR CLAIMS_MANAGER WHERE CLAIM_ID = :claim_id UPDATE CLAIMS SET STATUS = 'CANCELED'"
$ws.Range("F19").Characters(71, 37).Font.Size = 10
$ws.Range("F19").Characters(71, 37).Font.Name = "Arial Unicode MS"

$ws.Range("F20").Value = "UPDATE CLAIMS_MANAGER SET STATUS = 'UNDER_INVESTIGATION' WHERE CLAIM_ID = :claim_id UPDATE CLAIMS SET STATUS = 'UNDER_INVESTIGATION'"
$ws.Range("F20").Characters(85, 48).Font.Size = 10
$ws.Range("F20").Characters(85, 48).Font.Name = "Arial Unicode MS"

$ws.Range("F21").Value = "UPDATE OR DELETE ... WHERE CLAIM_ID IN (...) INSERT INTO AUDIT TABLE"
$ws.Range("F21").Characters(46, 23).Font.Size = 10
$ws.Range("F21").Characters(46, 23).Font.Name = "Arial Unicode MS"

$ws.Range("F23").Value = "UPDATE CLAIMS SET STATUS = 'REJECTED' WHERE STATUS = 'OUTSTANDING' AND AGE > :threshold UPDATE CLAIMS_MANAGER SET STATUS = 'CLOSED' WHERE STATUS = 'UNDER_INVESTIGATION' AND AGE > :threshold"
$ws.Range("F23").Characters(89, 101).Font.Size = 10
$ws.Range("F23").Characters(89, 101).Font.Name = "Arial Unicode MS"

# ---- Selection / view ----
$ws.Range("E33").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100

# ---- Column widths (AutoFit based on content, like original "bestFit") ----
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
